# Actualización de logs y cifras del 4 de Julio 2022
# Updates the "Julio" sheet: fills in the day's figures for 2022-07-04
# (row 6, corresponding to Excel serial date 44746) and restores the
# sheet's active selection to L5.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Julio")

# Fecha (A6) is already populated with 44746 (2022-07-04); fill in the
# day's counts across the remaining columns.
$ws.Range("B6").Value2 = 97   # Total Eventos Identificados
$ws.Range("C6").Value2 = 3    # Volcaduras
$ws.Range("D6").Value2 = 3    # Peaton Atropellado
$ws.Range("E6").Value2 = 0    # Motocilista Atropellado
$ws.Range("F6").Value2 = 0    # Ciclista Atropellado
$ws.Range("G6").Value2 = 1    # Peaton Fallecido
$ws.Range("H6").Value2 = 0    # Ciclista Fallecido
$ws.Range("I6").Value2 = 0    # Motociclista Fallecido
$ws.Range("J6").Value2 = 0    # Automovilista Fallecido

# Move/restore the active selection on the "Julio" sheet as left by the author.
$ws.Activate()
$ws.Range("L5").Select()
